$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164478063583374
$ws.Range("B1").Value = 2.279306650161743
$ws.Range("C1").Value = 4.431197643280029
$ws.Range("D1").Value = 3.46958327293396
$ws.Range("E1").Value = 1.228003025054932
